# Apply the metrics_2_9 update: reorder model names in column A (rows 2-26)
# and set the new constant metric values (columns B-I) for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New order of model names for rows 2..26 (row -> model name)
$newNames = @{
    2  = "model_2_9_0"
    3  = "model_2_9_22"
    4  = "model_2_9_21"
    5  = "model_2_9_20"
    6  = "model_2_9_19"
    7  = "model_2_9_18"
    8  = "model_2_9_17"
    9  = "model_2_9_16"
    10 = "model_2_9_15"
    11 = "model_2_9_14"
    12 = "model_2_9_13"
    13 = "model_2_9_23"
    14 = "model_2_9_12"
    15 = "model_2_9_10"
    16 = "model_2_9_9"
    17 = "model_2_9_8"
    18 = "model_2_9_7"
    19 = "model_2_9_6"
    20 = "model_2_9_5"
    21 = "model_2_9_4"
    22 = "model_2_9_3"
    23 = "model_2_9_2"
    24 = "model_2_9_11"
    25 = "model_2_9_1"
    26 = "model_2_9_24"
}

# New constant metric values applied to every data row (B..I)
$valB = 0.09932080507882668
$valC = 0.03320776602868791
$valD = -0.117905973768957
$valE = -0.007641264704611572
$valF = 0.9967864155769348
$valG = 1.576687693595886
$valH = 0.7930145263671875
$valI = 1.207900285720825

for ($row = 2; $row -le 26; $row++) {
    $ws.Range("A$row").Value = $newNames[$row]
    $ws.Range("B$row").Value = $valB
    $ws.Range("C$row").Value = $valC
    $ws.Range("D$row").Value = $valD
    $ws.Range("E$row").Value = $valE
    $ws.Range("F$row").Value = $valF
    $ws.Range("G$row").Value = $valG
    $ws.Range("H$row").Value = $valH
    $ws.Range("I$row").Value = $valI
}
